$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in the title cell (row 1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 00:04"

# Update country data rows whose ranking position / daily figures changed.
# Data rows are sorted descending by total cases (column B); a handful of
# countries changed rank (Guinea-Bisau, Gabon, Liberia, Santo Tome y Principe,
# Uganda moved up) while others simply got refreshed case counts.

$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1289028
$ws.Cells.Item(4, 3).Value = 25936
$ws.Cells.Item(4, 4).Value = 215580
$ws.Cells.Item(4, 5).Value = 996657
$ws.Cells.Item(4, 6).Value = 17009
$ws.Cells.Item(4, 7).Value = 1992
$ws.Cells.Item(4, 8).Value = 76791

$ws.Cells.Item(10, 1).Value = "Alemania"
$ws.Cells.Item(10, 2).Value = 169430
$ws.Cells.Item(10, 3).Value = 1268
$ws.Cells.Item(10, 4).Value = 139900
$ws.Cells.Item(10, 5).Value = 22138
$ws.Cells.Item(10, 6).Value = 1823
$ws.Cells.Item(10, 7).Value = 117
$ws.Cells.Item(10, 8).Value = 7392

$ws.Cells.Item(12, 1).Value = "Brasil"
$ws.Cells.Item(12, 2).Value = 132367
$ws.Cells.Item(12, 3).Value = 5756
$ws.Cells.Item(12, 4).Value = 51370
$ws.Cells.Item(12, 5).Value = 71943
$ws.Cells.Item(12, 6).Value = 8318
$ws.Cells.Item(12, 7).Value = 466
$ws.Cells.Item(12, 8).Value = 9054

$ws.Cells.Item(35, 1).Value = "Japon"
$ws.Cells.Item(35, 2).Value = 15477
$ws.Cells.Item(35, 3).Value = 224
$ws.Cells.Item(35, 4).Value = 4918
$ws.Cells.Item(35, 5).Value = 9982
$ws.Cells.Item(35, 6).Value = 308
$ws.Cells.Item(35, 7).Value = 21
$ws.Cells.Item(35, 8).Value = 577

$ws.Cells.Item(108, 1).Value = "Burkina Faso"
$ws.Cells.Item(108, 2).Value = 736
$ws.Cells.Item(108, 3).Value = 7
$ws.Cells.Item(108, 4).Value = 562
$ws.Cells.Item(108, 5).Value = 126
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 48

$ws.Cells.Item(117, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(117, 2).Value = 564
$ws.Cells.Item(117, 3).Value = 89
$ws.Cells.Item(117, 4).Value = 25
$ws.Cells.Item(117, 5).Value = 537
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0
$ws.Cells.Item(117, 8).Value = 2

$ws.Cells.Item(118, 1).Value = "Gabon"
$ws.Cells.Item(118, 2).Value = 504
$ws.Cells.Item(118, 3).Value = 107
$ws.Cells.Item(118, 4).Value = 110
$ws.Cells.Item(118, 5).Value = 386
$ws.Cells.Item(118, 6).Value = 1
$ws.Cells.Item(118, 7).Value = 2
$ws.Cells.Item(118, 8).Value = 8

$ws.Cells.Item(119, 1).Value = "Jordania"
$ws.Cells.Item(119, 2).Value = 494
$ws.Cells.Item(119, 3).Value = 21
$ws.Cells.Item(119, 4).Value = 381
$ws.Cells.Item(119, 5).Value = 104
$ws.Cells.Item(119, 6).Value = 5
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 9

$ws.Cells.Item(120, 1).Value = "Malta"
$ws.Cells.Item(120, 2).Value = 486
$ws.Cells.Item(120, 3).Value = 2
$ws.Cells.Item(120, 4).Value = 413
$ws.Cells.Item(120, 5).Value = 68
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 5

$ws.Cells.Item(121, 1).Value = "Tanzania"
$ws.Cells.Item(121, 2).Value = 480
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 167
$ws.Cells.Item(121, 5).Value = 297
$ws.Cells.Item(121, 6).Value = 7
$ws.Cells.Item(121, 7).Value = 0
$ws.Cells.Item(121, 8).Value = 16

$ws.Cells.Item(122, 1).Value = "Jamaica"
$ws.Cells.Item(122, 2).Value = 478
$ws.Cells.Item(122, 3).Value = 5
$ws.Cells.Item(122, 4).Value = 57
$ws.Cells.Item(122, 5).Value = 412
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 9

$ws.Cells.Item(123, 1).Value = "Paraguay"
$ws.Cells.Item(123, 2).Value = 462
$ws.Cells.Item(123, 3).Value = 22
$ws.Cells.Item(123, 4).Value = 148
$ws.Cells.Item(123, 5).Value = 304
$ws.Cells.Item(123, 6).Value = 9
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 10

$ws.Cells.Item(124, 1).Value = "Tayikistan"
$ws.Cells.Item(124, 2).Value = 461
$ws.Cells.Item(124, 3).Value = 82
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 449
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 4
$ws.Cells.Item(124, 8).Value = 12

$ws.Cells.Item(125, 1).Value = "Taiwan"
$ws.Cells.Item(125, 2).Value = 440
$ws.Cells.Item(125, 3).Value = 1
$ws.Cells.Item(125, 4).Value = 347
$ws.Cells.Item(125, 5).Value = 87
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 6

$ws.Cells.Item(126, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(126, 2).Value = 439
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(126, 4).Value = 13
$ws.Cells.Item(126, 5).Value = 422
$ws.Cells.Item(126, 6).Value = 0
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 4

$ws.Cells.Item(127, 1).Value = "Reunion"
$ws.Cells.Item(127, 2).Value = 427
$ws.Cells.Item(127, 3).Value = 2
$ws.Cells.Item(127, 4).Value = 354
$ws.Cells.Item(127, 5).Value = 73
$ws.Cells.Item(127, 6).Value = 3
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

$ws.Cells.Item(131, 1).Value = "Isla de Man"
$ws.Cells.Item(131, 2).Value = 329
$ws.Cells.Item(131, 3).Value = 2
$ws.Cells.Item(131, 4).Value = 271
$ws.Cells.Item(131, 5).Value = 35
$ws.Cells.Item(131, 6).Value = 19
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 23

$ws.Cells.Item(134, 1).Value = "Ruanda"
$ws.Cells.Item(134, 2).Value = 271
$ws.Cells.Item(134, 3).Value = 3
$ws.Cells.Item(134, 4).Value = 133
$ws.Cells.Item(134, 5).Value = 138
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 0

$ws.Cells.Item(141, 1).Value = "Liberia"
$ws.Cells.Item(141, 2).Value = 189
$ws.Cells.Item(141, 3).Value = 11
$ws.Cells.Item(141, 4).Value = 79
$ws.Cells.Item(141, 5).Value = 90
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 20

$ws.Cells.Item(142, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(142, 2).Value = 187
$ws.Cells.Item(142, 3).Value = 13
$ws.Cells.Item(142, 4).Value = 4
$ws.Cells.Item(142, 5).Value = 179
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = 4

$ws.Cells.Item(143, 1).Value = "Islas Feroe"
$ws.Cells.Item(143, 2).Value = 187
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 185
$ws.Cells.Item(143, 5).Value = 2
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 0

$ws.Cells.Item(144, 1).Value = "Martinica"
$ws.Cells.Item(144, 2).Value = 183
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(144, 4).Value = 83
$ws.Cells.Item(144, 5).Value = 86
$ws.Cells.Item(144, 6).Value = 3
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 14

$ws.Cells.Item(145, 1).Value = "Birmania"
$ws.Cells.Item(145, 2).Value = 176
$ws.Cells.Item(145, 3).Value = 15
$ws.Cells.Item(145, 4).Value = 62
$ws.Cells.Item(145, 5).Value = 108
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 6

$ws.Cells.Item(149, 1).Value = "Gibraltar"
$ws.Cells.Item(149, 2).Value = 144
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 141
$ws.Cells.Item(149, 5).Value = 3
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 0

$ws.Cells.Item(159, 1).Value = "Uganda"
$ws.Cells.Item(159, 2).Value = 101
$ws.Cells.Item(159, 3).Value = 1
$ws.Cells.Item(159, 4).Value = 55
$ws.Cells.Item(159, 5).Value = 46
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 0

$ws.Cells.Item(160, 1).Value = "Aruba"
$ws.Cells.Item(160, 2).Value = 101
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 89
$ws.Cells.Item(160, 5).Value = 9
$ws.Cells.Item(160, 6).Value = 4
$ws.Cells.Item(160, 7).Value = 1
$ws.Cells.Item(160, 8).Value = 3
